$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E49").Value = "DONE"

$ws.Range("C50").Value = 29.1
$ws.Range("E50").Value = "DONE"

$ws.Range("C51").Value = 29.2
$ws.Range("D51").Value = "Statistic filter"
$ws.Range("E51").Value = "DONE"
$ws.Range("F51").Value = "Frontend"
$ws.Range("F51").Style = $ws.Range("F50").Style

$ws.Range("C52").Value = 30
$ws.Range("D52").Value = "Create, Edit, Delete for Auth"
$ws.Range("F52").Value = "Frontend"

$ws.Range("C52:F52").Borders.LineStyle = 1
